$d = $word.ActiveDocument

# --- Edit 1: Append " (da completare)" after the Protocolli bullet text, with
#             "da completare" in bold red. ---
$rng = $d.Content
$found1 = $rng.Find.Execute("Aggiunta gestione Protocolli (Menù --> Protocolli)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the Protocolli bullet text"
}
$rng.Collapse(0)

$rng.InsertAfter(" (")
$rng.Collapse(0)

$rng.InsertAfter("da completare")
$rng.Bold = 1
$rng.Font.Color = 255

$rng.Collapse(0)
$rng.InsertAfter(")")

# --- Edit 2: Add a new bulleted paragraph after the IBAN bullet. ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Aggiunto campo IBAN alla tabella COLLABORATORI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the IBAN bullet text"
}
$rng2.Collapse(0)

$rng2.InsertParagraphAfter()
$rng2.Collapse(0)
$null = $rng2.Move(1, 1)

$rng2.InsertAfter("Aggiustamenti layout verticale offerta e appuntamento da agenda")
